# Segg the new transmittals and included the data files for poc
#
# Adds three new "Transmittals_New_*" test cases (Correspondence,
# ConsultantAdvice, ChangeNote) to both the DataFetchFlag and DataFetchXL
# sheets, right after the existing "Transmittals_New" row and before the
# "Documents_New" row (which shifts down).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("DataFetchFlag")
$ws2 = $wb.Worksheets.Item("DataFetchXL")

# ---------------------------------------------------------------------
# Make room: duplicate row 3 (Documents_New) into rows 3..5 on both
# sheets, pushing the original row 3 content down to row 6. This keeps
# cell formatting (style) consistent with the existing rows.
# ---------------------------------------------------------------------
$ws1.Range("A3:B3").Copy() | Out-Null
$ws1.Range("A3:B5").Insert(-4121) | Out-Null

# Remember the hyperlink cell style so it can be re-applied after the
# hyperlinks collection is rebuilt below.
$hlStyle = $ws2.Range("B2").Style

$ws2.Range("A3:C3").Copy() | Out-Null
$ws2.Range("A3:C5").Insert(-4121) | Out-Null

# ---------------------------------------------------------------------
# Fill in the new test case rows.
# ---------------------------------------------------------------------
$ws1.Range("A3").Value = "Transmittals_New_Correspondence"
$ws1.Range("A4").Value = "Transmittals_New_ConsultantAdvice"

$ws2.Range("A3").Value = "Transmittals_New_Correspondence"
$ws2.Range("B3").Value = '\\src\\com\\proj\\suiteTRANSMITTALS\\testdata\\TransmittalsTestData-Correspondence.xlsx'
$ws2.Range("A4").Value = "Transmittals_New_ConsultantAdvice"
$ws2.Range("B4").Value = '\\src\\com\\proj\\suiteTRANSMITTALS\\testdata\\TransmittalsTestData-ConsultantAdvice.xlsx'

$ws1.Range("A5").Value = "Transmittals_New_ChangeNote"
$ws2.Range("A5").Value = "Transmittals_New_ChangeNote"
$ws2.Range("B5").Value = '\\src\\com\\proj\\suiteTRANSMITTALS\\testdata\\TransmittalsTestData-ChangeNote.xlsx'

$ws1.Range("B3").Value = "XL"
$ws1.Range("B4").Value = "XL"
$ws1.Range("B5").Value = "XL"

$ws2.Range("C3").Value = "Transmittals_New"
$ws2.Range("C4").Value = "Transmittals_New"
$ws2.Range("C5").Value = "Transmittals_New"

# ---------------------------------------------------------------------
# The hyperlink that used to live on B3 (Documents_New) needs to end up
# on B6 now that the row shifted down; the engine doesn't move hyperlink
# anchors on insert, so rebuild the whole collection in final order.
# ---------------------------------------------------------------------
$ws2.Hyperlinks.Delete() | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("B2"), 'file:///\\src\com\proj\suiteTRANSMITTALS\testdata\TransmittalsTestData.xlsx') | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("B6"), 'file:///\\src\com\proj\suiteDOCS\testdata\DocumentsTestData.xlsx') | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("B3"), 'file:///\\src\com\proj\suiteTRANSMITTALS\testdata\TransmittalsTestData-Correspondence.xlsx') | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("B4"), 'file:///\\src\com\proj\suiteTRANSMITTALS\testdata\TransmittalsTestData-ConsultantAdvice.xlsx') | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("B5"), 'file:///\\src\com\proj\suiteTRANSMITTALS\testdata\TransmittalsTestData-ChangeNote.xlsx') | Out-Null

$ws2.Range("B2").Style = $hlStyle
$ws2.Range("B3").Style = $hlStyle
$ws2.Range("B4").Style = $hlStyle
$ws2.Range("B5").Style = $hlStyle
$ws2.Range("B6").Style = $hlStyle

$ws2.Range("B12").Select() | Out-Null

# Restore the active sheet/selection to DataFetchFlag!A5 (it was the
# active sheet before this edit and stays active afterwards).
$ws1.Activate() | Out-Null
$ws1.Range("A5").Select() | Out-Null
